# Audi.xlsx edit: rename the sheet, drop the bold header styling, and
# move the active selection, per the target commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Sheet1" -> "CarDetails"
$ws.Name = "CarDetails"

# Header cells (A1:C1) lose their bold font.
$ws.Range("A1:C1").Font.Bold = $false

# Selection moves from B4 to A5.
[void]$ws.Range("A5").Select()
